$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B18").Value = "MEC-1NA-M.T.F."
$ws.Range("F18").Value = "-"
$ws.Range("B19").Value = "MEC-1NA-M.T.F."
$ws.Range("B20").Value = "-"
